$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column S: add a new "2022" column, mirroring the existing R column ---

# Row 3 (thin/empty divider row before the bottom border) - copy format from R3
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)   # xlPasteFormats

# Row 4 (year headers). R4 changes from the "20xx" look to match D4:P4, then S4 = 2022
$ws.Range("D4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Row 5 (first data row). R5 now matches Q5's style, S5 = 1.8
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 1.8

# Row 6 (second data row). R6 now matches L6:P6's style, S6 = 8.4
$ws.Range("L6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 8.4

$excel.CutCopyMode = $false

$ws.Range("S3").Select()
